# Update cryptos list figures (prices + 1h volume %) per latest scrape,
# and swap the Binance-PegBSC-USD / Fetch.AI rows (30 <-> 31).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.583.78"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.65%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.523.10"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.33%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.95%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "152.62"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.82%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.537"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.525.48"
$ws.Range("D9").Style = "Normal"
$ws.Range("E10").Value = "  +0.43%  "
$ws.Range("E11").Value = "  -1.38%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.28"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.12%  "
$ws.Range("E13").Value = "  -1.21%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "29.36"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.00%  "
$ws.Range("E15").Value = "  +1.18%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.990.51"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.67%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.222.68"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.24%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.528.86"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.41%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.82"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.55%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.96"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.31%  "
$ws.Range("E21").Value = "  +3.13%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "328.89"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.52%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.23"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.96"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.75%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "65.41"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.18%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "648.93"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.21%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0000104"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.52%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.648.50"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.19%  "
$ws.Range("B30").Value = "Fetch.AI"
$ws.Range("C30").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.50"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.86%  "
$ws.Range("B31").Value = "Binance-PegBSC-USD"
$ws.Range("C31").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.05%  "
$ws.Range("E32").Value = "  +0.40%  "
$ws.Range("E33").Value = "  +1.46%  "
$ws.Range("E34").Value = "  +2.16%  "
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("E36").Value = "  +1.23%  "
$ws.Range("E37").Value = "  +1.26%  "
$ws.Range("E38").Value = "  +1.83%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "154.06"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.75%  "
$ws.Range("E40").Value = "  +0.76%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "18.92"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.11%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.82"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.35%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.79"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.47%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "163.25"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.71%  "
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₆0301"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.82%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "15.53"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.92%  "
$ws.Range("E48").Value = "  +1.73%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "21.37"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.83%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.618"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.74%  "
$ws.Range("E51").Value = "  +1.26%  "
